$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Set Runmode column (C) to "Y" for all rows except row 14 (already "Y")
for ($r = 2; $r -le 18; $r++) {
    if ($r -ne 14) {
        $ws.Range("C$r").Value = "Y"
    }
}

# Update selection to reflect the new active range C2:C18
$ws.Range("C2:C18").Select()
